# Updated cryptos list on Mon Feb 26 02:51:51 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    # Force the literal text value (avoids Excel auto-converting
    # numeric-looking strings like "384.67" into real numbers), then
    # strip the "quote prefix" style Excel applies so the cell keeps
    # its original (default) style.
    $c = $ws.Range($addr)
    $c.Value = "'" + $val
    $c.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextCell "D2" "51.493.07"
$ws.Range("E2").Value = "  -0.27%  "

# Row 3 - Ethereum
Set-TextCell "D3" "3.101.20"
$ws.Range("E3").Value = "  +2.40%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.02%  "

# Row 5 - BNB
Set-TextCell "D5" "384.67"
$ws.Range("E5").Value = "  +1.15%  "

# Row 6 - Solana
Set-TextCell "D6" "103.14"
$ws.Range("E6").Value = "  -0.26%  "

# Row 7 - XRP
Set-TextCell "D7" "0.541"
$ws.Range("E7").Value = "  -0.90%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.04%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  -2.28%  "

# Row 10 - Avalanche
Set-TextCell "D10" "37.22"
$ws.Range("E10").Value = "  +0.97%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  -0.18%  "

# Row 12 - Dogecoin
Set-TextCell "D12" "0.0854"
$ws.Range("E12").Value = "  -0.63%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-TextCell "D13" "3.592.54"
$ws.Range("E13").Value = "  +2.48%  "

# Row 14 - Chainlink
Set-TextCell "D14" "18.63"
$ws.Range("E14").Value = "  +0.27%  "

# Row 15 - Polkadot
$ws.Range("E15").Value = "  +0.56%  "

# Row 16 - now WrappedEther (was Uniswap)
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextCell "D16" "3.097.89"
$ws.Range("E16").Value = "  +1.28%  "

# Row 17 - now Uniswap (was WrappedEther)
$ws.Range("B17").Value = "Uniswap"
$ws.Range("C17").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextCell "D17" "11.36"
$ws.Range("E17").Value = "  +6.33%  "

# Row 18 - Polygon
Set-TextCell "D18" "0.995"
$ws.Range("E18").Value = "  +0.03%  "

# Row 19 - WrappedBTC
Set-TextCell "D19" "51.521.41"
$ws.Range("E19").Value = "  -0.39%  "

# Row 20 - ImmutableX
Set-TextCell "D20" "3.31"
$ws.Range("E20").Value = "  +7.85%  "

# Row 21 - ShibaInu
$ws.Range("E21").Value = "  -0.12%  "

# Row 22 - InternetComputer(DFINITY)
Set-TextCell "D22" "12.33"
$ws.Range("E22").Value = "  -1.52%  "

# Row 23 - Litecoin
Set-TextCell "D23" "69.94"
$ws.Range("E23").Value = "  -0.48%  "

# Row 24 - BitcoinCash
Set-TextCell "D24" "265.81"
$ws.Range("E24").Value = "  -1.15%  "

# Row 25 - PancakeSwap
Set-TextCell "D25" "3.11"
$ws.Range("E25").Value = "  -2.15%  "

# Row 26 - Filecoin
Set-TextCell "D26" "8.10"
$ws.Range("E26").Value = "  -3.26%  "

# Row 27 - EthereumClassic
Set-TextCell "D27" "27.03"
$ws.Range("E27").Value = "  +3.03%  "

# Row 28 - RenderToken
Set-TextCell "D28" "7.24"
$ws.Range("E28").Value = "  -3.15%  "

# Row 29 - Dai
$ws.Range("E29").Value = "  +0.09%  "

# Row 30 - Kaspa
$ws.Range("E30").Value = "  -2.18%  "

# Row 31 - Hedera
$ws.Range("E31").Value = "  -2.61%  "

# Row 32 - Cosmos
$ws.Range("E32").Value = "  -0.24%  "

# Row 33 - InjectiveProtocol
Set-TextCell "D33" "35.31"
$ws.Range("E33").Value = "  +2.81%  "

# Row 34 - VeChain
Set-TextCell "D34" "0.0471"
$ws.Range("E34").Value = "  +4.32%  "

# Row 35 - Toncoin
$ws.Range("E35").Value = "  +0.50%  "

# Row 36 - OKB
Set-TextCell "D36" "50.35"
$ws.Range("E36").Value = "  -1.57%  "

# Row 37 - FirstDigitalUSD
$ws.Range("E37").Value = "  -0.20%  "

# Row 38 - LidoDAOToken
Set-TextCell "D38" "3.37"
$ws.Range("E38").Value = "  +1.84%  "

# Row 39 - TheGraph
Set-TextCell "D39" "0.298"
$ws.Range("E39").Value = "  +5.54%  "

# Row 40 - ARBITRUM
$ws.Range("E40").Value = "  +1.02%  "

# Row 41 - Monero
Set-TextCell "D41" "129.07"
$ws.Range("E41").Value = "  +1.67%  "

# Row 42 - Stellar
$ws.Range("E42").Value = "  -0.69%  "

# Row 43 - Celestia
Set-TextCell "D43" "16.55"
$ws.Range("E43").Value = "  -4.40%  "

# Row 44 - Stacks
$ws.Range("E44").Value = "  -3.51%  "

# Row 45 - NEARProtocol
Set-TextCell "D45" "3.68"
$ws.Range("E45").Value = "  -2.67%  "

# Row 46 - EnergySwap
Set-TextCell "D46" "22.27"
$ws.Range("E46").Value = "  +1.85%  "

# Row 47 - ApeXProtocol
$ws.Range("E47").Value = "  +3.16%  "

# Row 48 - WEMIXToken
$ws.Range("E48").Value = "  +1.71%  "

# Row 49 - Maker
Set-TextCell "D49" "2.056.61"
$ws.Range("E49").Value = "  +0.72%  "

# Row 50 - RocketPoolETH
Set-TextCell "D50" "3.412.56"
$ws.Range("E50").Value = "  +2.26%  "

# Row 51 - BEAM
Set-TextCell "D51" "0.0327"
$ws.Range("E51").Value = "  +1.06%  "
